# Insert a new "Title and Content" slide before the current last slide
# (so the slide order becomes: 1, 2, 3, NEW, old-4), reusing a table laid
# out inside the content placeholder's area, mirroring how the author
# dropped a table into slide 4 of the deck (pushing the old chart slide
# down to position 5).

$p = $ppt.ActivePresentation

# 2 = ppLayoutObject -> "Title and Content" custom layout.
$s = $p.Slides.Add(4, 2)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Text = "Slide title"
$titleRange.Font.Name = "Calibri"
$titleRange.Font.NameFarEast = "Calibri"
$titleRange.Font.NameComplexScript = "Calibri"

# --- Replace the content placeholder with a table ---------------------
$s.Shapes.Item(2).Delete()

# Dimensions match the table inserted into the content placeholder area
# (in points; AddTable takes points, 914400 EMU = 72 pt -> EMU/12700).
$tbl = $s.Shapes.AddTable(2, 3, 66, 143.75, 640, 58.4)
$table = $tbl.Table

$table.Cell(1,1).Shape.TextFrame.TextRange.Text = "Test table"
$table.Cell(1,2).Shape.TextFrame.TextRange.Text = "Names"
$table.Cell(1,3).Shape.TextFrame.TextRange.Text = "Emails"

$table.Cell(2,2).Shape.TextFrame.TextRange.Text = "{{ program.users.name }}"
$table.Cell(2,3).Shape.TextFrame.TextRange.Text = "{{ program.users.email }}"
